$wb = $excel.ActiveWorkbook

# --- 1. Germany sheet: selection was a whole-column pick (A1:XFD1048576),
#        narrow it down to the sheet's real used range (A1:D12) ---
$germany = $wb.Worksheets.Item("Germany")
$germany.Range("A1:D12").Select()

# --- 2. Add the new "Portugal" market sheet. It is a copy of the existing
#        "Swiss" sheet (same layout/styles/merged cells), placed after the
#        last tab, then re-labelled with Portugal's own data. ---
$swiss = $wb.Worksheets.Item("Swiss")
$swiss.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$portugal = $wb.Worksheets.Item($wb.Worksheets.Count)
$portugal.Name = "Portugal"

# Market name + Jira/ticket reference for the Portugal market
$portugal.Range("B2").Value = "Portugal Market"
$portugal.Range("B4").Value = "NGC-3479/T2459"

# Narrower columns than the Swiss sheet (tuned to fit the Portugal content)
$portugal.Columns.Item(1).ColumnWidth = 13.166666666666668
$portugal.Columns.Item(2).ColumnWidth = 19.833333333333336
$portugal.Columns.Item(3).ColumnWidth = 12.333333333333332
$portugal.Columns.Item(4).ColumnWidth = 14.0

# Rows 3-5 grew taller (wrapped content) compared to the Swiss sheet
$portugal.Rows.Item(3).RowHeight = 28.8
$portugal.Rows.Item(4).RowHeight = 28.8
$portugal.Rows.Item(5).RowHeight = 28.8

# Leave the cursor on B4:B5, and make Portugal the active/visible tab
$portugal.Range("B4:B5").Select()
$portugal.Activate()
